$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "actual" (column B) values for several rows
$ws.Range("B5").Value = 318
$ws.Range("B6").Value = 328
$ws.Range("B17").Value = 339
$ws.Range("B18").Value = 240
$ws.Range("B28").Value = 271
$ws.Range("B29").Value = 405
$ws.Range("B40").Value = 462
$ws.Range("B41").Value = 390
$ws.Range("B52").Value = 391
$ws.Range("B53").Value = 273

# Update prediction columns (C: prediksi_mean, D: lower, E: upper) for rows 62-85
$ws.Range("C62").Value = 200
$ws.Range("D62").Value = 119.3230632560374
$ws.Range("E62").Value = 280.6769367439626
$ws.Range("C63").Value = 233
$ws.Range("D63").Value = 152.3230632560374
$ws.Range("E63").Value = 313.6769367439626
$ws.Range("C64").Value = 391
$ws.Range("D64").Value = 310.3230632560374
$ws.Range("E64").Value = 471.6769367439626
$ws.Range("C65").Value = 273
$ws.Range("D65").Value = 192.3230632560374
$ws.Range("E65").Value = 353.6769367439626
$ws.Range("C66").Value = 225
$ws.Range("D66").Value = 144.3230632560374
$ws.Range("E66").Value = 305.6769367439626
$ws.Range("C67").Value = 239
$ws.Range("D67").Value = 158.3230632560374
$ws.Range("E67").Value = 319.6769367439626
$ws.Range("C68").Value = 244
$ws.Range("D68").Value = 163.3230632560374
$ws.Range("E68").Value = 324.6769367439626
$ws.Range("C69").Value = 206
$ws.Range("D69").Value = 125.3230632560374
$ws.Range("E69").Value = 286.6769367439626
$ws.Range("C70").Value = 202
$ws.Range("D70").Value = 121.3230632560374
$ws.Range("E70").Value = 282.6769367439626
$ws.Range("C71").Value = 209
$ws.Range("D71").Value = 128.3230632560374
$ws.Range("E71").Value = 289.6769367439626
$ws.Range("C72").Value = 225
$ws.Range("D72").Value = 144.3230632560374
$ws.Range("E72").Value = 305.6769367439626
$ws.Range("C73").Value = 222
$ws.Range("D73").Value = 141.3230632560374
$ws.Range("E73").Value = 302.6769367439626
$ws.Range("C74").Value = 200
$ws.Range("D74").Value = 85.90558188597133
$ws.Range("E74").Value = 314.0944181140287
$ws.Range("C75").Value = 233
$ws.Range("D75").Value = 118.9055818859713
$ws.Range("E75").Value = 347.0944181140287
$ws.Range("C76").Value = 391
$ws.Range("D76").Value = 276.9055818859713
$ws.Range("E76").Value = 505.0944181140287
$ws.Range("C77").Value = 273
$ws.Range("D77").Value = 158.9055818859713
$ws.Range("E77").Value = 387.0944181140287
$ws.Range("C78").Value = 225
$ws.Range("D78").Value = 110.9055818859713
$ws.Range("E78").Value = 339.0944181140287
$ws.Range("C79").Value = 239
$ws.Range("D79").Value = 124.9055818859713
$ws.Range("E79").Value = 353.0944181140287
$ws.Range("C80").Value = 244
$ws.Range("D80").Value = 129.9055818859713
$ws.Range("E80").Value = 358.0944181140287
$ws.Range("C81").Value = 206
$ws.Range("D81").Value = 91.90558188597133
$ws.Range("E81").Value = 320.0944181140287
$ws.Range("C82").Value = 202
$ws.Range("D82").Value = 87.90558188597133
$ws.Range("E82").Value = 316.0944181140287
$ws.Range("C83").Value = 209
$ws.Range("D83").Value = 94.90558188597133
$ws.Range("E83").Value = 323.0944181140287
$ws.Range("C84").Value = 225
$ws.Range("D84").Value = 110.9055818859713
$ws.Range("E84").Value = 339.0944181140287
$ws.Range("C85").Value = 222
$ws.Range("D85").Value = 107.9055818859713
$ws.Range("E85").Value = 336.0944181140287
